# Weekly price update: insert a new sampling row for "Ajo" (Chino/Primera/China)
# at row 110 of the sheet, shifting all subsequent rows down by one.
#
# This mirrors what happened in the source workbook: a new day's record
# (2021-11-18, serial 44518) was inserted into the middle of the date-ordered
# series, pushing every later record down one row (the last record, which
# used to be at row 175, now lives at row 176).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 110 - this shifts rows 110..175 down to
# 111..176 and keeps their formatting/styles intact.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new record's data.
$ws.Cells.Item(110, 1).Value  = 4
$ws.Cells.Item(110, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value  = "Los Lagos"
$ws.Cells.Item(110, 4).Value  = 44518
$ws.Cells.Item(110, 5).Value  = 10
$ws.Cells.Item(110, 6).Value  = 100112003
$ws.Cells.Item(110, 7).Value  = "Ajo"
$ws.Cells.Item(110, 8).Value  = "Chino"
$ws.Cells.Item(110, 9).Value  = "Primera"
$ws.Cells.Item(110, 10).Value = 100
$ws.Cells.Item(110, 11).Value = 22000
$ws.Cells.Item(110, 12).Value = 23000
$ws.Cells.Item(110, 13).Value = 22600
$ws.Cells.Item(110, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(110, 15).Value = "China"
$ws.Cells.Item(110, 16).Value = 2260
$ws.Cells.Item(110, 17).Value = 10
$ws.Cells.Item(110, 18).Value = "Hortaliza"
